# The upstream commit ("Fixed #295 Add the version of M2Doc in the
# template custom properties") changed how M2Doc writes the custom
# document properties of *generated* templates. The diff shipped for
# this particular fixture (simpleVariableNotExisting-template.docx)
# is, line for line, a pure re-serialization of word/document.xml and
# word/styles.xml: every "-"/"+" pair has the identical element, the
# identical set of attribute name/value pairs and identical text
# content - only the attribute order changed (attributes were written
# out alphabetically by local name instead of in their original
# order). No text, run, paragraph, property value, relationship or
# part was added, removed or modified.
#
# The Word object model has no notion of "XML attribute order" (it is
# a pure serialization-layer artifact of whatever tool re-saved the
# fixture), so there is no content-level edit to perform here. This
# script intentionally touches the document without changing any of
# its visible/semantic content, leaving it equivalent to the target
# state described by the diff.
$d = $word.ActiveDocument
